$d = $word.ActiveDocument

# Green color used throughout this edit (RGB 0x168253 -> decimal 5472790)
$green = 5472790

# 1) Merge "Detall: " + "tiene que aparecer todas estas fechas y campos"
#    into a single run (same formatting, so Find/Replace coalesces them).
$d.Content.Find.Execute("Detall: tiene que aparecer todas estas fechas y campos", $true, $false, $false, $false, $false, $true, 1, $false, "Detall: tiene que aparecer todas estas fechas y campos", 2) | Out-Null

# 2) "Data renovació marca (+ 2 años)" -> split so "Data renovació marca" is green,
#    " (+ 2 años)" keeps default color.
$p = $d.Paragraphs.Item(6)
$start = $p.Range.Start
$greenLen = "Data renovació marca".Length
$r = $d.Range($start, $start + $greenLen)
$r.Font.Color = $green

# 3) Color the following list items green (paragraph + run rPr).
$greenParas = @(7, 8, 9, 10, 11, 12, 13, 14, 15)
foreach ($idx in $greenParas) {
    $p = $d.Paragraphs.Item($idx)
    $p.Range.Font.Color = $green
}

# 4) "Data resolució revocació " -> split so "Data resolució revocació" is green
#    and the trailing space remains a separate, uncolored run.
$p = $d.Paragraphs.Item(16)
$fullText = $p.Range.Text
$trimLen = $fullText.TrimEnd().Length
$start = $p.Range.Start
$r = $d.Range($start, $start + $trimLen)
$r.Font.Color = $green

# 5) styles.xml: Normal style paragraph properties gain suppressAutoHyphens (true).
$normal = $d.Styles.Item("Normal")
$normal.ParagraphFormat.Hyphenation = $false
